# Updated the release locations

$wb = $excel.ActiveWorkbook

# --- spawning-sites: bump release number (D) and particle count (E) for every site row ---
$wsSpawn = $wb.Worksheets.Item("spawning-sites")
$wsSpawn.Range("D2:D18").Value = 5
$wsSpawn.Range("E2:E18").Value = 4000

# --- basic-biology: new species + citations ---
$wsBio = $wb.Worksheets.Item("basic-biology")
$wsBio.Range("B2").Value = "Abudefduf vaigiensis "
$wsBio.Range("C2").Value = "Wellington & Victor 1989"
$wsBio.Range("B3").Value = 18.3
$wsBio.Range("B4").Value = 1.5
$wsBio.Range("C8").Value = "Murphy et al. (2007)"
$wsBio.Range("B10").Value = "n/a"

# --- basic-physical: turbulence-horizontal value change ---
$wsPhys = $wb.Worksheets.Item("basic-physical")
$wsPhys.Range("B2").Value = 0.3

# --- restore/update the selections (active cell) recorded per-sheet ---
$wsBio.Range("B11").Select() | Out-Null
$wsPhys.Range("E7").Select() | Out-Null

# spawning-sites stays the active tab/sheet, so select it last
$wsSpawn.Range("L17").Select() | Out-Null
